$wb = $excel.ActiveWorkbook

# Update the migration period label on the "parametrosInicio" sheet (row 10)
# from "MIGRACIONES SGV DICIEMBRE 2022 28.12.2022" to "MIGRACIONES SGV ENERO 2023 10.01.2023"
$wsParams = $wb.Worksheets.Item("parametrosInicio")
$wsParams.Range("B10").Value = "MIGRACIONES SGV ENERO 2023 10.01.2023"

# Restore/keep the selection on "Rutas" sheet at B3 first (it stops being the active sheet)
$wsRutas = $wb.Worksheets.Item("Rutas")
$wsRutas.Range("B3").Select()

# Select C10 on "parametrosInicio" and make it the active sheet (select/activate last so
# it remains the sheet that is active/selected when the workbook is saved)
$wsParams.Range("C10").Select()
$wsParams.Activate()
